$wb = $excel.ActiveWorkbook

# The two worksheets ("positions" and "positionsC") hold the same CV table.
# Rows 28-30 list three publications; this change re-orders them so the
# (corrected) 2013 PNAS paper becomes the first of the three, pushing the
# other two rows down by one position.
$sheetNames = @("positions", "positionsC")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Capture the current ("before") values of the three affected rows.
    $row28 = @{
        C = $ws.Range("C28").Value()
        D = $ws.Range("D28").Value()
        F = $ws.Range("F28").Value()
        G = $ws.Range("G28").Value()
        H = $ws.Range("H28").Value()
        J = $ws.Range("J28").Value()
    }
    $row29 = @{
        C = $ws.Range("C29").Value()
        D = $ws.Range("D29").Value()
        F = $ws.Range("F29").Value()
        G = $ws.Range("G29").Value()
        H = $ws.Range("H29").Value()
        J = $ws.Range("J29").Value()
    }
    $row30 = @{
        C = $ws.Range("C30").Value()
        D = $ws.Range("D30").Value()
        F = $ws.Range("F30").Value()
        G = $ws.Range("G30").Value()
        H = $ws.Range("H30").Value()
        J = $ws.Range("J30").Value()
    }

    # New row 28 = old row 30, but with the publication year corrected
    # from 2011 to 2013 (the citation's actual publication year).
    $ws.Range("C28").Value = $row30.C
    $ws.Range("D28").Value = $row30.D
    $ws.Range("F28").Value = 2013
    $ws.Range("G28").Value = 2013
    $ws.Range("H28").Value = $row30.H
    $ws.Range("J28").Value = $row30.J

    # New row 29 = old row 28.
    $ws.Range("C29").Value = $row28.C
    $ws.Range("D29").Value = $row28.D
    $ws.Range("F29").Value = $row28.F
    $ws.Range("G29").Value = $row28.G
    $ws.Range("H29").Value = $row28.H
    $ws.Range("J29").Value = $row28.J

    # New row 30 = old row 29.
    $ws.Range("C30").Value = $row29.C
    $ws.Range("D30").Value = $row29.D
    $ws.Range("F30").Value = $row29.F
    $ws.Range("G30").Value = $row29.G
    $ws.Range("H30").Value = $row29.H
    $ws.Range("J30").Value = $row29.J
}

# Update sheet view / selection state to match the saved workbook:
#  - "positionsC" keeps a plain (non-active) view with the whole of row 28 selected
#  - "positions" becomes the active/selected tab, with F28 selected
$wsPositions = $wb.Worksheets.Item("positions")
$wsPositionsC = $wb.Worksheets.Item("positionsC")

$wsPositionsC.Activate()
$wsPositionsC.Rows.Item(28).Select()

$wsPositions.Activate()
$wsPositions.Range("F28").Select()
